# Remove the "ICPetcdHighNumberOfFailedGRPCRequests" alert rows (warning +
# critical variants) because of false alerts on ICP.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 29 and 30 hold the ICPetcdHighNumberOfFailedGRPCRequests alert
# (warning severity, then critical severity). Delete both, shifting the
# remaining rows (ICPetcdHighNumberOfFailedProposals, ICPetcdHighFsyncDurations,
# ICPetcdHighCommitDurations) up by two.
$deleteRange = $ws.Range("A29:E30")
$deleteRange.EntireRow.Delete()

# Mirror the selection state recorded after the deletion.
$ws.Range("A29:XFD30").Select()
